$d = $word.ActiveDocument

# --- Remove the stray "_GoBack" bookmark left over in the body (an empty,
#     center-aligned paragraph right before the ROI formula). Word normally
#     drops this automatically; here we do it explicitly. ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Fix the paging: the section used to carry a single "default" footer
#     (with a stale "start at 58" page numbering, left over from when this
#     chapter was a standalone doc). Splitting into even / default / first
#     page headers & footers and resetting the starting page number to 1
#     resolves the paging glitch. ---
$sec = $d.Sections(1)

$sec.Headers(1).Range.Text = ""
$sec.Headers(2).Range.Text = ""
$sec.Headers(3).Range.Text = ""
$sec.Footers(1).Range.Text = ""
$sec.Footers(2).Range.Text = ""
$sec.Footers(3).Range.Text = ""

$sec.Headers(1).Range.Paragraphs(1).Style = $d.Styles("header")
$sec.Headers(2).Range.Paragraphs(1).Style = $d.Styles("header")
$sec.Headers(3).Range.Paragraphs(1).Style = $d.Styles("header")
$sec.Footers(1).Range.Paragraphs(1).Style = $d.Styles("footer")
$sec.Footers(2).Range.Paragraphs(1).Style = $d.Styles("footer")
$sec.Footers(3).Range.Paragraphs(1).Style = $d.Styles("footer")

# Re-plant the "_GoBack" bookmark (Word re-creates this on every save to
# mark the last edit position) inside the primary footer, now that it is
# an empty paragraph.
$sec.Footers(1).Range.Bookmarks.Add("_GoBack") | Out-Null

$sec.Headers(1).PageNumbers.StartingNumber = 1
